$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stage the new certificate number / date as genuine text (via a formula that
# returns a string) in a couple of scratch cells far away from the table, so
# that the leading zero of the cert number and the dotted date are preserved
# as plain text instead of being auto-converted to a number/date.
$ws.Range("ZZ1").Formula = '="0898840"'
$ws.Range("ZZ2").Formula = '="2020.07.05"'

# Insert a new row at row 5 (shifts existing rows 5-20 down to 6-21)
$ws.Rows.Item(5).Insert()

# Copy the formatting of row 4 (a "company certificate" row) into the new row 5
$ws.Range("A4:D4").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)

# Bring in the staged text values (values only, keeping row 5's formatting)
$ws.Range("ZZ1").Copy()
$ws.Range("A5").PasteSpecial(-4163)
$ws.Range("ZZ2").Copy()
$ws.Range("B5").PasteSpecial(-4163)

$ws.Range("C5").Value = "Git"
$ws.Range("D5").Value = "Geekbrains"

# Clean up the scratch cells
$ws.Range("ZZ1:ZZ2").Clear()

# Restore the active cell selection as recorded after the edit
$ws.Range("C4").Select()
